# Updates cryptocurrency price/volume data cells on Sheet1 (Price column D, Volume(1h) column E).
# Values are written as text (matching the existing inline-string cell type) via the leading
# apostrophe text-qualifier, then the cell style is reset to "Normal" so no stray number-format
# / quote-prefix styling is left behind on cells that previously had none.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.994.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.61%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.640.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.67%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.44%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'215.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.77%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.19%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.40%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.41%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.48%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.869.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.70%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.642.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.75%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.543"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = "'  +1.10%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.21%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.086.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.94%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'194.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.65%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.54%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.16%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.67%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +4.96%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.11%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.35%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'142.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.19%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.60%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.72%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.75%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0496"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.32%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.23%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.04%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.19%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +1.20%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.36%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.128.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.67%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.538"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.51%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.23%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.27%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.68%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'98.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.03%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.796"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.55%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.777.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.65%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +4.08%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'56.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.27%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'1.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.31%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.95%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.13%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.56%  "
$ws.Range("E51").Style = "Normal"
